# Update weekly excess mortality analysis (CBS "Berekening oversterfte" sheet)
# - revise several already-published weekly "Waargenomen" (observed) figures
#   in column G (CBS's normal weekly data revisions)
# - add week 49 as a new data row (F/G/H/I), which pushes the totals row
#   down by one
# - move the cell selection to J41, matching where the user ended up after
#   adding the new row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 41 for the new "week 49" entry. This pushes the old
# row 41 (blank spacer) and row 42 (totals) down to 42 and 43 respectively.
$ws.Rows("41").Insert()

# --- Revisions to previously reported weekly "Waargenomen" (column G) ---
$ws.Range("G8").Value  = 4305
$ws.Range("G21").Value = 2528
$ws.Range("G23").Value = 2667
$ws.Range("G24").Value = 2640
$ws.Range("G26").Value = 2853
$ws.Range("G31").Value = 2891
$ws.Range("G33").Value = 3019
$ws.Range("G34").Value = 3212
$ws.Range("G35").Value = 3444
$ws.Range("G36").Value = 3674
$ws.Range("G37").Value = 3587
$ws.Range("G38").Value = 3552
$ws.Range("G39").Value = 3315
$ws.Range("G40").Value = 3373

# --- New week 49 row ---
$ws.Range("F41").Value = 49
$ws.Range("G41").Value = 3448
$ws.Range("H41").Value = 3037
$ws.Range("I41").Formula = "=G41-H41"

# Move the selection to where the editor left off (one column to the right
# of the new "Oversterfte" figure).
$ws.Range("J41").Select() | Out-Null
